$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 35, shifting existing rows 35..110 down to 36..111
$ws.Rows.Item(35).Insert()

# Populate the newly inserted row 35 with the new data record
$ws.Cells.Item(35, 1).Value2  = 3
$ws.Cells.Item(35, 2).Value2  = "Femacal de La Calera"
$ws.Cells.Item(35, 3).Value2  = "Coquimbo"
$ws.Cells.Item(35, 4).Value2  = 44622
$ws.Cells.Item(35, 5).Value2  = 5
$ws.Cells.Item(35, 6).Value2  = "Fruta"
$ws.Cells.Item(35, 7).Value2  = 100107
$ws.Cells.Item(35, 8).Value2  = "Otros"
$ws.Cells.Item(35, 9).Value2  = 100107011
$ws.Cells.Item(35, 10).Value2 = "Tuna"
$ws.Cells.Item(35, 11).Value2 = "Sin especificar"
$ws.Cells.Item(35, 12).Value2 = "Primera"
$ws.Cells.Item(35, 13).Value2 = 87
$ws.Cells.Item(35, 14).Value2 = 14000
$ws.Cells.Item(35, 15).Value2 = 14000
$ws.Cells.Item(35, 16).Value2 = 14000
$ws.Cells.Item(35, 17).Value2 = "$/caja 16 kilos"
$ws.Cells.Item(35, 18).Value2 = "Cabildo"
$ws.Cells.Item(35, 19).Value2 = 875
$ws.Cells.Item(35, 20).Value2 = 16
